$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "DS 인력풀이 문제? 교육이 문제? 나라가 문제?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/ds-manpower-education-country/#utm_source=rss&utm_medium=rss&utm_campaign=ds-manpower-education-country"

$ws.Range("D51").Value = "[MySQL] 테이블 복원 중에 만난 ERROR 3546 해결 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/MySQL-%ED%85%8C%EC%9D%B4%EB%B8%94-%EB%B3%B5%EC%9B%90-%EC%A4%91%EC%97%90-%EB%A7%8C%EB%82%9C-ERROR-3546-%ED%95%B4%EA%B2%B0-%EB%B0%A9%EB%B2%95"
